$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card8")

# Insert a new column before column N (14th column), shifting old N ("Correction ") to O
$ws.Columns.Item(14).Insert()

# New N1 header (the newly inserted column)
$ws.Cells.Item(1, 14).Value = "Correction"

# Fill the new N column (rows 2-13) with "nan" to match data pattern
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}
